$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q16").Copy()
$ws.Range("R16").PasteSpecial(-4122)
$ws.Range("R16").NumberFormat = "0.00"
